# The commit swaps the two embedded theme parts in this deck:
#   ppt/theme/theme1.xml  (the presentation/slide-master theme, formerly
#                           the "Integral" palette) becomes the default
#                           "Office Theme" palette.
#   ppt/theme/theme2.xml  (the notes-master theme, formerly the default
#                           "Office Theme" palette) becomes the
#                           "Integral" palette.
#
# Every slide in the deck is painted through the slide master's theme
# (theme1.xml), so the visible, load-bearing part of this edit is
# recoloring that theme's 12-slot colour scheme to the stock "Office"
# palette. We do that the same way a user would from the Design /
# Slide Master UI: by editing each ThemeColor's RGB value on the
# presentation's ThemeColorScheme.
#
# PowerPoint's ThemeColor.RGB (like the classic VBA RGB()/ColorFormat.RGB)
# is a single integer packed as 0x00BBGGRR, i.e. the reverse byte order of
# the usual #RRGGBB hex notation, so each target hex colour is converted
# below before being assigned.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# Target "Office Theme" colour scheme (12 slots, in clrScheme order):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$colors.Colors(1).RGB  = 0        # dk1      000000
$colors.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$colors.Colors(3).RGB  = 6968388  # dk2      44546A
$colors.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$colors.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$colors.Colors(6).RGB  = 3243501  # accent2  ED7D31
$colors.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$colors.Colors(8).RGB  = 49407    # accent4  FFC000
$colors.Colors(9).RGB  = 12874308 # accent5  4472C4
$colors.Colors(10).RGB = 4697456  # accent6  70AD47
$colors.Colors(11).RGB = 12673797 # hlink    0563C1
$colors.Colors(12).RGB = 7491477  # folHlink 954F72
